$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Coin"/"Link"/"Price"/"Volume(1h)" cells (refreshed crypto price feed).
# Price/Volume columns are stored as literal text in this workbook, so force
# text formatting ("@") before writing each numeric-looking value -- otherwise
# Excel auto-converts a typed "285.49" / "-10.56%" into a number.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '285.49'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '-10.56%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '40.10'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '-3.17%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.041'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '-3.84%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.07267'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '4.304'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '-0.43%'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.506'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '-10.73%'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.9130'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '-3.91%'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.1200'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '-3.10%'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1700'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-7.29%'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08561'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '-7.22%'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.04171'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-4.54%'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.1050'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '-0.09%'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.001274'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '-0.27%'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.006008'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '0.88%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.401'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '1.85%'
$ws.Range('B17').Value = 'BTSEToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.397'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '-1.16%'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.3261'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '-2.94%'
$ws.Range('B19').Value = 'MCDex'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.863'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '2.07%'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.1355'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '0.22%'
$ws.Range('B21').Value = 'ZBToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.2890'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '2.47%'
$ws.Range('B22').Value = 'CoinExToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.03851'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-4.43%'
$ws.Range('B23').Value = 'BitKan'
$ws.Range('C23').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.001271'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '0.50%'
$ws.Range('B24').Value = 'HotbitToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.003783'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '-8.07%'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0001284'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '1.16%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0003734'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02275'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '-10.56%'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.04912'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '-8.27%'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.007061'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '255.02%'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.007674'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '-1.02%'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1262'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '-4.33%'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.007394'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '0.69%'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.006944'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '-8.43%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.3083'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '-10.44%'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00006398'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '-4.58%'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.00000000752'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '0.35%'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '20.83%'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.00002106'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '0.35%'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0002006'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '0.35%'
